{"js": "// Update the date label and every \"a\u00f7b=c, d\" answer cell in the table.\n// Each (oldText -> newText) pair below corresponds 1:1, in document order,\n// to the runs touched by the diff.\nconst replacements = [\n  [\"2023-09-06 Wednesday\", \"2023-09-07 Thursday\"],\n  [\"33\u00f76=5, 3\", \"58\u00f77=8, 2\"],\n  [\"26\u00f76=4, 2\", \"14\u00f79=1, 5\"],\n  [\"10\u00f76=1, 4\", \"82\u00f77=11, 5\"],\n  [\"99\u00f72=49, 1\", \"20\u00f73=6, 2\"],\n  [\"92\u00f79=10, 2\", \"87\u00f79=9, 6\"],\n  [\"64\u00f78=8, 0\", \"33\u00f75=6, 3\"],\n  [\"37\u00f76=6, 1\", \"14\u00f72=7, 0\"],\n  [\"99\u00f73=33, 0\", \"81\u00f74=20, 1\"],\n  [\"29\u00f77=4, 1\", \"69\u00f75=13, 4\"],\n  [\"95\u00f74=23, 3\", \"30\u00f79=3, 3\"],\n  [\"79\u00f74=19, 3\", \"90\u00f76=15, 0\"],\n  [\"99\u00f74=24, 3\", \"80\u00f77=11, 3\"],\n  [\"48\u00f79=5, 3\", \"86\u00f72=43, 0\"],\n  [\"11\u00f79=1, 2\", \"43\u00f74=10, 3\"],\n  [\"25\u00f79=2, 7\", \"79\u00f79=8, 7\"],\n  [\"42\u00f78=5, 2\", \"71\u00f79=7, 8\"],\n  [\"14\u00f76=2, 2\", \"19\u00f73=6, 1\"],\n  [\"11\u00f78=1, 3\", \"81\u00f77=11, 4\"],\n  [\"28\u00f75=5, 3\", \"71\u00f75=14, 1\"],\n  [\"46\u00f75=9, 1\", \"95\u00f79=10, 5\"],\n  [\"47\u00f78=5, 7\", \"55\u00f77=7, 6\"],\n  [\"22\u00f79=2, 4\", \"77\u00f77=11, 0\"],\n  [\"87\u00f79=9, 6\", \"89\u00f75=17, 4\"],\n  [\"21\u00f74=5, 1\", \"73\u00f74=18, 1\"],\n  [\"58\u00f74=14, 2\", \"57\u00f72=28, 1\"],\n];\n\nconst body = context.document.body;\n\n// Resolve every search range up front (against the still-untouched\n// document), so an earlier replacement's new text can never be picked\n// up as the match for a later search (this matters because some new\n// values equal other rows' old values, e.g. \"87\u00f79=9, 6\").\nconst searchResults = replacements.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((r) => r.load(\"items\"));\nawait context.sync();\n\nsearchResults.forEach((r, i) => {\n  const [oldText, newText] = replacements[i];\n  if (r.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for ${JSON.stringify(oldText)}, found ${r.items.length}`\n    );\n  }\n  r.items[0].insertText(newText, Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Update the date label (first paragraph) and every \"a\u00f7b=c, d\" answer in\n# the table. Cells are addressed by (row, column) rather than by searching\n# for their old text, because a couple of the new values equal another\n# cell's *old* value (e.g. \"87\u00f79=9, 6\"); a naive sequential Find/Replace\n# over the whole document could re-match text that an earlier replacement\n# had just inserted. Writing straight to each Cell.Range.Text sidesteps\n# that entirely and leaves the run/paragraph formatting untouched.\n\n$d = $word.ActiveDocument\n\n# First paragraph: the date line.\n$d.Paragraphs.Item(1).Range.Text = \"2023-09-07 Thursday\"\n\n$tbl = $d.Tables.Item(1)\n\n# The 25 new answers, in reading order (row by row, left to right) \u2014 only\n# the rows that actually hold an answer (every 4th row is a blank spacer).\n$newValues = @(\n    \"58\u00f77=8, 2\", \"14\u00f79=1, 5\", \"82\u00f77=11, 5\", \"20\u00f73=6, 2\", \"87\u00f79=9, 6\",\n    \"33\u00f75=6, 3\", \"14\u00f72=7, 0\", \"81\u00f74=20, 1\", \"69\u00f75=13, 4\", \"30\u00f79=3, 3\",\n    \"90\u00f76=15, 0\", \"80\u00f77=11, 3\", \"86\u00f72=43, 0\", \"43\u00f74=10, 3\", \"79\u00f79=8, 7\",\n    \"71\u00f79=7, 8\", \"19\u00f73=6, 1\", \"81\u00f77=11, 4\", \"71\u00f75=14, 1\", \"95\u00f79=10, 5\",\n    \"55\u00f77=7, 6\", \"77\u00f77=11, 0\", \"89\u00f75=17, 4\", \"73\u00f74=18, 1\", \"57\u00f72=28, 1\"\n)\n\n# Find the rows that already contain text (the blank spacer rows stay empty).\n$dataRowIndexes = @()\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    $firstCellText = $tbl.Cell($r, 1).Range.Text.TrimEnd([char]7, [char]13, [char]10)\n    if ($firstCellText.Length -gt 0) {\n        $dataRowIndexes += $r\n    }\n}\n\n$colCount = $tbl.Columns.Count\n$valueIndex = 0\nforeach ($row in $dataRowIndexes) {\n    for ($col = 1; $col -le $colCount; $col++) {\n        $tbl.Cell($row, $col).Range.Text = $newValues[$valueIndex]\n        $valueIndex++\n    }\n}\n"}
